$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the spicule type label in row 2 from "small tylostyle" to "Oxea"
$ws.Range("A2").Value = "Oxea"
